$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Manchester City v Leeds United" fixture row (row 1).
$ws.Rows(1).Delete()

# Remove the "Arsenal v Manchester City" fixture row
# (was row 3, now row 2 after the previous delete shifted rows up).
$ws.Rows(2).Delete()

# The Brentford fixture kickoff time moved from 20:00 to 19:45.
# (was row 5, now row 3 after the two deletes above).
$ws.Range("B3").Value = "09/02/2022 19:45 | Premier League"

Write-Host "Done updating fixture list."
